$wb = $excel.ActiveWorkbook

# --- PIR sheet: add rows 627-640 ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A627:F640").NumberFormat = "@"
$ws.Range("A627:A640").Value = '2026-02-06'
$ws.Range("C627:C640").Value = '10:00'
$ws.Range("D627:D640").Value = 'Bathroom'
$ws.Cells.Item(627,2).Value = '10:33:30'
$ws.Cells.Item(628,2).Value = '10:33:35'
$ws.Cells.Item(629,2).Value = '10:33:36'
$ws.Cells.Item(630,2).Value = '10:33:40'
$ws.Cells.Item(631,2).Value = '10:33:45'
$ws.Cells.Item(632,2).Value = '10:33:50'
$ws.Cells.Item(633,2).Value = '10:33:55'
$ws.Cells.Item(634,2).Value = '10:34:00'
$ws.Cells.Item(635,2).Value = '10:34:05'
$ws.Cells.Item(636,2).Value = '10:34:09'
$ws.Cells.Item(637,2).Value = '10:34:15'
$ws.Cells.Item(638,2).Value = '10:34:19'
$ws.Cells.Item(639,2).Value = '10:34:24'
$ws.Cells.Item(640,2).Value = '10:34:26'
$ws.Cells.Item(627,5).Value = 'No Motion'
$ws.Cells.Item(628,5).Value = 'No Motion'
$ws.Cells.Item(629,5).Value = 'No Motion'
$ws.Cells.Item(630,5).Value = 'No Motion'
$ws.Cells.Item(631,5).Value = 'No Motion'
$ws.Cells.Item(632,5).Value = 'No Motion'
$ws.Cells.Item(633,5).Value = 'No Motion'
$ws.Cells.Item(634,5).Value = 'No Motion'
$ws.Cells.Item(635,5).Value = 'No Motion'
$ws.Cells.Item(636,5).Value = 'Motion Detected'
$ws.Cells.Item(637,5).Value = 'No Motion'
$ws.Cells.Item(638,5).Value = 'Motion Detected'
$ws.Cells.Item(639,5).Value = 'No Motion'
$ws.Cells.Item(640,5).Value = 'Motion Detected'
$ws.Cells.Item(627,6).Value = 'Inactive'
$ws.Cells.Item(628,6).Value = 'Inactive'
$ws.Cells.Item(629,6).Value = 'Inactive'
$ws.Cells.Item(630,6).Value = 'Inactive'
$ws.Cells.Item(631,6).Value = 'Inactive'
$ws.Cells.Item(632,6).Value = 'Inactive'
$ws.Cells.Item(633,6).Value = 'Inactive'
$ws.Cells.Item(634,6).Value = 'Inactive'
$ws.Cells.Item(635,6).Value = 'Inactive'
$ws.Cells.Item(636,6).Value = 'Active'
$ws.Cells.Item(637,6).Value = 'Inactive'
$ws.Cells.Item(638,6).Value = 'Active'
$ws.Cells.Item(639,6).Value = 'Inactive'
$ws.Cells.Item(640,6).Value = 'Active'

# --- Humidity sheet: add rows 450-459 ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A450:F459").NumberFormat = "@"
$ws.Range("A450:A459").Value = '2026-02-06'
$ws.Range("C450:C459").Value = '10:00'
$ws.Range("D450:D459").Value = 'Bathroom'
$ws.Cells.Item(450,2).Value = '10:33:32'
$ws.Cells.Item(451,2).Value = '10:33:38'
$ws.Cells.Item(452,2).Value = '10:33:42'
$ws.Cells.Item(453,2).Value = '10:33:46'
$ws.Cells.Item(454,2).Value = '10:33:51'
$ws.Cells.Item(455,2).Value = '10:34:07'
$ws.Cells.Item(456,2).Value = '10:34:11'
$ws.Cells.Item(457,2).Value = '10:34:16'
$ws.Cells.Item(458,2).Value = '10:34:21'
$ws.Cells.Item(459,2).Value = '10:34:27'
$ws.Cells.Item(450,5).Value = '67.2%'
$ws.Cells.Item(451,5).Value = '65.6%'
$ws.Cells.Item(452,5).Value = '67.1%'
$ws.Cells.Item(453,5).Value = '67.2%'
$ws.Cells.Item(454,5).Value = '67.1%'
$ws.Cells.Item(455,5).Value = '66.1%'
$ws.Cells.Item(456,5).Value = '67.3%'
$ws.Cells.Item(457,5).Value = '66.4%'
$ws.Cells.Item(458,5).Value = '67.7%'
$ws.Cells.Item(459,5).Value = '68.0%'
$ws.Cells.Item(450,6).Value = 'Active'
$ws.Cells.Item(451,6).Value = 'Active'
$ws.Cells.Item(452,6).Value = 'Active'
$ws.Cells.Item(453,6).Value = 'Active'
$ws.Cells.Item(454,6).Value = 'Active'
$ws.Cells.Item(455,6).Value = 'Active'
$ws.Cells.Item(456,6).Value = 'Active'
$ws.Cells.Item(457,6).Value = 'Active'
$ws.Cells.Item(458,6).Value = 'Active'
$ws.Cells.Item(459,6).Value = 'Active'

# --- Temperature sheet: add rows 449-458 ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A449:F458").NumberFormat = "@"
$ws.Range("A449:A458").Value = '2026-02-06'
$ws.Range("C449:C458").Value = '10:00'
$ws.Range("D449:D458").Value = 'Bathroom'
$ws.Cells.Item(449,2).Value = '10:33:33'
$ws.Cells.Item(450,2).Value = '10:33:39'
$ws.Cells.Item(451,2).Value = '10:33:43'
$ws.Cells.Item(452,2).Value = '10:33:48'
$ws.Cells.Item(453,2).Value = '10:33:53'
$ws.Cells.Item(454,2).Value = '10:34:08'
$ws.Cells.Item(455,2).Value = '10:34:12'
$ws.Cells.Item(456,2).Value = '10:34:18'
$ws.Cells.Item(457,2).Value = '10:34:22'
$ws.Cells.Item(458,2).Value = '10:34:28'
$ws.Cells.Item(449,5).Value = '28.5C'
$ws.Cells.Item(450,5).Value = '28.5C'
$ws.Cells.Item(451,5).Value = '28.5C'
$ws.Cells.Item(452,5).Value = '28.6C'
$ws.Cells.Item(453,5).Value = '28.5C'
$ws.Cells.Item(454,5).Value = '28.5C'
$ws.Cells.Item(455,5).Value = '28.5C'
$ws.Cells.Item(456,5).Value = '28.4C'
$ws.Cells.Item(457,5).Value = '28.5C'
$ws.Cells.Item(458,5).Value = '28.5C'
$ws.Cells.Item(449,6).Value = 'Active'
$ws.Cells.Item(450,6).Value = 'Active'
$ws.Cells.Item(451,6).Value = 'Active'
$ws.Cells.Item(452,6).Value = 'Active'
$ws.Cells.Item(453,6).Value = 'Active'
$ws.Cells.Item(454,6).Value = 'Active'
$ws.Cells.Item(455,6).Value = 'Active'
$ws.Cells.Item(456,6).Value = 'Active'
$ws.Cells.Item(457,6).Value = 'Active'
$ws.Cells.Item(458,6).Value = 'Active'
